$d = $word.ActiveDocument

# --- Step 1: "eprésentée par Monsieur " -> "eprésentée par " ---
$r1 = $d.Content
$null = $r1.Find.Execute("eprésentée par Monsieur ", $true, $false, $false, $false, $false,
                          $true, 1, $false, "", 0)
$rSeg1 = $d.Range($r1.Start, $r1.End)
$rSeg1.Text = "eprésentée par "

# --- Step 2: insert "Madame " (non-bold) right before the director's bold name ---
$r2 = $d.Content
$null = $r2.Find.Execute("El Hadji Mamadou FAYE", $true, $false, $false, $false, $false,
                          $true, 1, $false, "", 0)
$insPoint = $d.Range($r2.Start, $r2.Start)
$insPoint.InsertBefore("Madame ")

# --- Step 3: replace the bold name, extending the bold run to swallow the
#             trailing ", " so "Jenny MVOU, " stays bold together ---
$r3 = $d.Content
$null = $r3.Find.Execute("El Hadji Mamadou FAYE", $true, $false, $false, $false, $false,
                          $true, 1, $false, "", 0)
$null = $r3.MoveEnd(1, 2)
$r3.Text = "Jenny MVOU, "

# --- Step 4: update the role description to the feminine wording ---
$r4 = $d.Content
$null = $r4.Find.Execute("en qualité de Directeur Général, dument habilité aux fins des présentes",
                          $true, $false, $false, $false, $false,
                          $true, 1, $false, "", 0)
$rSeg4 = $d.Range($r4.Start, $r4.End)
$rSeg4.Text = "en qualité de Directeur Général Adjointe, dument habilitée aux fins des présentes"

Write-Output "edit applied"
